# Claim VAT refund DIY invoice template - add 2 new columns
# (Total invoice amount swaps ahead of VAT amount, and a brand new
#  "Is the Invoice in your name? Y/N" column is appended.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Re-order / add the header row text.
#    Old layout: F=VAT amount, G=Total invoice amount (excluding VAT)
#    New layout: F=Total invoice amount (excluding VAT), G=VAT amount,
#                H=Is the Invoice in your name? Y/N
# ---------------------------------------------------------------------
$ws.Range("F1").Value = "Total invoice amount (excluding VAT)"
$ws.Range("G1").Value = "VAT amount"
$ws.Range("H1").Value = "Is the Invoice in your name? Y/N"

# ---------------------------------------------------------------------
# 2. Borders - give every header cell A1:G1 a thin box border
#    (bottom edge is already thin from the original template style,
#    so only top/left/right need to be added).
# ---------------------------------------------------------------------

# H1 first - left/top border added (bottom already thin); no right edge
# (right edge continues visually as I1's left edge). Doing this before
# A1:G1 lets the left+top+bottom combination be reused below.
$h1 = $ws.Range("H1")
$h1.Borders.Item(8).LineStyle = 1
$h1.Borders.Item(7).LineStyle = 1

$headerCols = @("A1","B1","C1","D1","E1","F1","G1")
foreach ($addr in $headerCols) {
    $rng = $ws.Range($addr)
    $rng.Borders.Item(8).LineStyle = 1   # top
    $rng.Borders.Item(7).LineStyle = 1   # left
    $rng.Borders.Item(10).LineStyle = 1  # right
}

# I1 - left border only; drop the inherited bottom border
$i1 = $ws.Range("I1")
$i1.Borders.Item(9).LineStyle = -4142
$i1.Borders.Item(7).LineStyle = 1

# J1:Z1 - no border at all (remove the inherited bottom border)
$ws.Range("J1:Z1").Borders.Item(9).LineStyle = -4142

# ---------------------------------------------------------------------
# 3. Column widths - column F widens slightly.
# ---------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 13.67

# ---------------------------------------------------------------------
# 4. Sheet view / selection / page setup / footer
# ---------------------------------------------------------------------
$ws.Range("J7").Select() | Out-Null

$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
$ps.CenterFooter = '&1#&"Calibri"&10&K000000OFFICIAL'
